$d = $word.ActiveDocument

$replacements = @(
    @{old = "72×56="; new = "42×54="},
    @{old = "54×31="; new = "11×78="},
    @{old = "39×34="; new = "16×13="},
    @{old = "41×88="; new = "22×88="},
    @{old = "60×53="; new = "41×60="},
    @{old = "98×13="; new = "78×14="},
    @{old = "82×29="; new = "90×24="},
    @{old = "91×91="; new = "57×26="},
    @{old = "80×15="; new = "37×69="},
    @{old = "65×69="; new = "96×69="},
    @{old = "35×97="; new = "54×41="},
    @{old = "15×41="; new = "32×37="},
    @{old = "54×30="; new = "97×97="},
    @{old = "66×33="; new = "92×27="},
    @{old = "73×33="; new = "90×96="},
    @{old = "94×69="; new = "11×71="},
    @{old = "39×70="; new = "87×31="},
    @{old = "16×49="; new = "52×70="},
    @{old = "49×69="; new = "70×11="},
    @{old = "86×34="; new = "53×78="},
    @{old = "77×55="; new = "85×37="},
    @{old = "93×37="; new = "40×59="},
    @{old = "53×68="; new = "38×25="},
    @{old = "13×46="; new = "19×46="},
    @{old = "31×90="; new = "91×74="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
